$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gUSD")

# ---------------------------------------------------------------------------
# New daily data rows (26-32) appended below the existing table (row 25 was
# the last populated row: A3:K25).
# ---------------------------------------------------------------------------

# Row 26
$ws.Range("B26").Value = 11.09
$ws.Range("C26").Value = 36.72
$ws.Range("D26").Value = 13.33
$ws.Range("E26").Value = 4.38
$ws.Range("F26").Value = 11.34

# Row 27
$ws.Range("B27").Value = 10.42
$ws.Range("C27").Value = 37.37
$ws.Range("D27").Value = 13.45
$ws.Range("E27").Value = 12.15
$ws.Range("F27").Value = 10.04

# Row 28
$ws.Range("B28").Value = 9.66
$ws.Range("C28").Value = 37.73
$ws.Range("D28").Value = 13.57
$ws.Range("E28").Value = 5.2
$ws.Range("F28").Value = 9.57

# Row 29
$ws.Range("B29").Value = 8.96
$ws.Range("C29").Value = 38.09
$ws.Range("D29").Value = 13.69
$ws.Range("E29").Value = 6.32
$ws.Range("F29").Value = 9.73

# Row 30
$ws.Range("B30").Value = 8.2
$ws.Range("C30").Value = 38.4
$ws.Range("D30").Value = 13.79
$ws.Range("E30").Value = 5.03
$ws.Range("F30").Value = 9.64

# Row 31
$ws.Range("B31").Value = 7.36
$ws.Range("C31").Value = 38.94
$ws.Range("D31").Value = 13.79
$ws.Range("E31").Value = 8.62
$ws.Range("F31").Value = 7.21

# Row 32
$ws.Range("B32").Value = 6.96
$ws.Range("C32").Value = 39.56
$ws.Range("D32").Value = 14.64
$ws.Range("E32").Value = 11.11
$ws.Range("F32").Value = 7.34

# Column A style (same as the rest of the column) for the new rows
$ws.Range("A26:A32").Style = $ws.Range("A25").Style

# ---------------------------------------------------------------------------
# Formulas.
# Row 26 continues the existing fill from row 25 one cell at a time (as if
# it had been typed/filled individually).
# ---------------------------------------------------------------------------
$ws.Range("A26").Formula = "=A25+1"
$ws.Range("H26").Formula = "=B26-B25"
$ws.Range("I26").Formula = "=C26-C25"
$ws.Range("K26").Formula = "=B26+C26"

# Rows 27-32 are then filled down together in one pass.
$ws.Range("A27:A32").Formula = "=A26+1"
$ws.Range("H27:H32").Formula = "=B27-B26"
$ws.Range("I27:I32").Formula = "=C27-C26"
$ws.Range("K27:K32").Formula = "=B27+C27"

# ---------------------------------------------------------------------------
# View: selection ends up on L32, and the sheet is scrolled back so that
# topLeftCell reverts to the default (A1-relative, i.e. no explicit scroll).
# ---------------------------------------------------------------------------
$ws.Range("L32").Select()
